$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.696222543716431
$ws.Range("B1").Value = 6.668671131134033
$ws.Range("C1").Value = 6.246330261230469
$ws.Range("D1").Value = 10.19025230407715
$ws.Range("E1").Value = 5.710779666900635
